$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Periodo Mora" (period) labels and their associated
# "Valor Mora" (value) figures for rows 16-20 (periods 2402-2406),
# swapping row 16 <-> row 20 and row 17 <-> row 19 (row 18 / 2404 stays put).

$ws.Range("E16").Value2 = "2406"
$ws.Range("F16").Value2 = 20800

$ws.Range("E17").Value2 = "2405"
$ws.Range("F17").Value2 = 52000

$ws.Range("E19").Value2 = "2403"
$ws.Range("F19").Value2 = 46400

$ws.Range("E20").Value2 = "2402"
$ws.Range("F20").Value2 = 46400
